# Clean up the vaccine name / brand labels across every sheet:
#  - collapse an embedded line break inside a cell into a single space
#    (e.g. "Hepatitis B [5]\nPediatric/Adolescent" / "Recombivax\nHB")
#  - strip the trailing footnote markers " [1]" .. " [5]" (the numbers in
#    brackets referred to footnotes that no longer exist on the sheet),
#    leaving the single space that preceded them
#
# Applied workbook-wide with Cells.Replace so it covers every table
# (Pediatric VFC Vaccine, Adult Vaccine, Pediatric Influenza Vaccine,
# Adult Influenza Vaccine) without having to hard-code cell addresses.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Turn embedded newlines within a cell into a plain space.
    $null = $ws.Cells.Replace("`n", " ")

    # Drop the " [1]" ... " [5]" footnote-reference suffixes.
    for ($n = 1; $n -le 5; $n++) {
        $marker = " [" + $n + "]"
        $null = $ws.Cells.Replace($marker, " ")
    }
}
